# "merged jacobs comments into the structure"
#
# Jacob's peer/self-assessment rows (row 3 in the "Self assesment" block,
# row 16 in the "Peer assessment" block) previously had empty Grade /
# Examples-of-actions cells. Fill them in with the grade "Insufficient"
# and the two comment variants that were merged in from Jacob's notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Peer  and self assessment"

# --- Row 3: Jacob, first ("Self assesment") block ---
$ws.Range("B3").Value = "Insufficient"
$ws.Range("C3").Value = "1) have not been part of the project since Riga. "

# --- Row 16: Jacob, second ("Peer assessment") block ---
$ws.Range("B16").Value = "Insufficient"
$ws.Range("C16").Value = "1) Has not been part of project since Riga."

# --- Update the saved selection/cursor position to reflect where the
#     editor ended up after entering the comments ---
$ws.Activate()
$ws.Range("A17").Select() | Out-Null
